# Apply cryptocurrency price/volume updates per the commit diff.
# Rows 32 and 33 also swap coin identity (Filecoin <-> ImmutableX).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "26.547.81"
Set-TextCell 2 5 "  -2.39%  "

# Row 3
Set-TextCell 3 4 "1.813.20"
Set-TextCell 3 5 "  -2.13%  "

# Row 4
Set-TextCell 4 4 "1.008"
Set-TextCell 4 5 "  +0.67%  "

# Row 5
Set-TextCell 5 4 "1.007"
Set-TextCell 5 5 "  +0.65%  "

# Row 6
Set-TextCell 6 4 "308.72"
Set-TextCell 6 5 "  -1.41%  "

# Row 7
Set-TextCell 7 4 "0.4576"
Set-TextCell 7 5 "  -1.34%  "

# Row 8
Set-TextCell 8 4 "0.3668"
Set-TextCell 8 5 "  -1.35%  "

# Row 9
Set-TextCell 9 4 "0.07157"
Set-TextCell 9 5 "  -1.65%  "

# Row 10
Set-TextCell 10 4 "0.8799"
Set-TextCell 10 5 "  -0.79%  "

# Row 11
Set-TextCell 11 4 "0.07799"
Set-TextCell 11 5 "  -0.43%  "

# Row 12
Set-TextCell 12 4 "19.40"
Set-TextCell 12 5 "  -3.09%  "

# Row 13
Set-TextCell 13 4 "1.761.90"
Set-TextCell 13 5 "  -0.70%  "

# Row 14
Set-TextCell 14 4 "5.293"
Set-TextCell 14 5 "  -1.54%  "

# Row 15
Set-TextCell 15 4 "6.383"
Set-TextCell 15 5 "  -2.05%  "

# Row 16
Set-TextCell 16 4 "86.24"
Set-TextCell 16 5 "  -5.13%  "

# Row 17
Set-TextCell 17 4 "1.008"
Set-TextCell 17 5 "  +0.66%  "

# Row 18
Set-TextCell 18 4 "0.000008595"
Set-TextCell 18 5 "  -3.55%  "

# Row 20
Set-TextCell 20 4 "26.562.87"
Set-TextCell 20 5 "  -2.43%  "

# Row 21
Set-TextCell 21 4 "14.29"
Set-TextCell 21 5 "  -2.90%  "

# Row 22
Set-TextCell 22 4 "5.011"
Set-TextCell 22 5 "  -1.07%  "

# Row 23
Set-TextCell 23 4 "10.46"
Set-TextCell 23 5 "  -0.33%  "

# Row 24
Set-TextCell 24 4 "1.982"
Set-TextCell 24 5 "  +1.59%  "

# Row 25
Set-TextCell 25 4 "150.95"
Set-TextCell 25 5 "  -0.49%  "

# Row 26
Set-TextCell 26 4 "18.01"
Set-TextCell 26 5 "  -1.97%  "

# Row 27
Set-TextCell 27 4 "2.071"
Set-TextCell 27 5 "  +1.63%  "

# Row 28
Set-TextCell 28 4 "112.67"
Set-TextCell 28 5 "  -2.65%  "

# Row 29
Set-TextCell 29 4 "4.863"
Set-TextCell 29 5 "  -3.79%  "

# Row 30
Set-TextCell 30 4 "0.08694"
Set-TextCell 30 5 "  -1.47%  "

# Row 31
Set-TextCell 31 4 "3.052"
Set-TextCell 31 5 "  -3.50%  "

# Row 32
Set-TextCell 32 2 "ImmutableX"
Set-TextCell 32 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell 32 4 "0.7350"
Set-TextCell 32 5 "  -4.04%  "

# Row 33
Set-TextCell 33 2 "Filecoin"
Set-TextCell 33 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 33 4 "4.480"
Set-TextCell 33 5 "  -0.45%  "

# Row 34
Set-TextCell 34 4 "1.119"
Set-TextCell 34 5 "  -4.04%  "

# Row 35
Set-TextCell 35 4 "1.005"
Set-TextCell 35 5 "  +0.64%  "

# Row 36
Set-TextCell 36 4 "2.561"
Set-TextCell 36 5 "  -6.00%  "

# Row 37
Set-TextCell 37 4 "1.082"
Set-TextCell 37 5 "  -2.24%  "

# Row 38
Set-TextCell 38 4 "0.01937"
Set-TextCell 38 5 "  -0.10%  "

# Row 39
Set-TextCell 39 4 "0.05116"
Set-TextCell 39 5 "  -1.65%  "

# Row 40
Set-TextCell 40 4 "2.897"
Set-TextCell 40 5 "  -1.34%  "

# Row 41
Set-TextCell 41 4 "6.985"
Set-TextCell 41 5 "  -0.37%  "

# Row 42
Set-TextCell 42 4 "0.5015"
Set-TextCell 42 5 "  -1.63%  "

# Row 43
Set-TextCell 43 4 "0.1565"
Set-TextCell 43 5 "  -3.78%  "

# Row 44
Set-TextCell 44 4 "8.156"
Set-TextCell 44 5 "  -2.99%  "

# Row 45
Set-TextCell 45 4 "1.008"
Set-TextCell 45 5 "  +0.81%  "

# Row 46
Set-TextCell 46 4 "0.4625"
Set-TextCell 46 5 "  -3.42%  "

# Row 47
Set-TextCell 47 4 "10.01"
Set-TextCell 47 5 "  -3.27%  "

# Row 48
Set-TextCell 48 4 "101.01"
Set-TextCell 48 5 "  -1.80%  "

# Row 49
Set-TextCell 49 4 "1.595"
Set-TextCell 49 5 "  -2.53%  "

# Row 50
Set-TextCell 50 4 "0.06017"
Set-TextCell 50 5 "  -3.10%  "

# Row 51
Set-TextCell 51 4 "64.14"
Set-TextCell 51 5 "  -2.23%  "
